# Update the "Förändrad" (Changed) date column (C) for every data row.
# All values in column C move from the Excel serial date 46060 (2026-02-07)
# to 46061 (2026-02-08), i.e. bump the "last changed" date by one day.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row on the sheet.
$lastRow = $ws.UsedRange.Rows.Count + $ws.UsedRange.Row - 1

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 46060) {
        $cell.Value2 = 46061
    }
}
